$wb = $excel.ActiveWorkbook

# --- Add the new worksheet "ODI Batting Extra" as the last sheet --------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Batting Extra"

# --- Header row -----------------------------------------------------------
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "BATTING_POSITION"
$ws.Range("C1").Value = "NUM_4"
$ws.Range("D1").Value = "NUM_6"
$ws.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Range("F1").Value = "MAN_OF_MATCH"

# Re-use the header formatting (bold, centered, bordered) already present
# on the other sheets in this workbook.
$headerSource = $wb.Worksheets.Item("ODI Bowling").Range("A1:F1")
$headerSource.Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

# --- Row 2 ------------------------------------------------------------
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "4592"
$ws.Range("A2").ClearFormats()

$ws.Range("B2").Value = 11

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "0"
$ws.Range("C2").ClearFormats()

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0"
$ws.Range("D2").ClearFormats()

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.49%"
$ws.Range("E2").ClearFormats()

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "NO"
$ws.Range("F2").ClearFormats()

# --- Row 3 ------------------------------------------------------------
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "4641"
$ws.Range("A3").ClearFormats()

$ws.Range("B3").Value = 11

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "0"
$ws.Range("C3").ClearFormats()

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "0"
$ws.Range("D3").ClearFormats()

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = ""
$ws.Range("E3").ClearFormats()

$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "NO"
$ws.Range("F3").ClearFormats()

# Restore the originally active sheet/selection.
$wb.Worksheets.Item("Player Info").Activate() | Out-Null
$wb.Worksheets.Item("Player Info").Range("A1").Select() | Out-Null
